$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows 716/717 column D uses the same date number format as D715
$ws.Range("D716").NumberFormat = $ws.Range("D715").NumberFormat
$ws.Range("D717").NumberFormat = $ws.Range("D715").NumberFormat

$ws.Cells.Item(667, 4).Value2 = 44783
$ws.Cells.Item(667, 14).Value2 = 20000
$ws.Cells.Item(667, 15).Value2 = 20000
$ws.Cells.Item(667, 16).Value2 = 20000
$ws.Cells.Item(667, 19).Value2 = 1667
$ws.Cells.Item(668, 4).Value2 = 44783
$ws.Cells.Item(668, 12).Value2 = 'Segunda'
$ws.Cells.Item(668, 13).Value2 = 108
$ws.Cells.Item(668, 14).Value2 = 20000
$ws.Cells.Item(668, 15).Value2 = 20000
$ws.Cells.Item(668, 16).Value2 = 20000
$ws.Cells.Item(668, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(668, 19).Value2 = 1429
$ws.Cells.Item(668, 20).Value2 = 14
$ws.Cells.Item(669, 4).Value2 = 44425
$ws.Cells.Item(669, 13).Value2 = 108
$ws.Cells.Item(670, 4).Value2 = 44377
$ws.Cells.Item(670, 12).Value2 = 'Primera'
$ws.Cells.Item(670, 13).Value2 = 162
$ws.Cells.Item(670, 14).Value2 = 17000
$ws.Cells.Item(670, 15).Value2 = 17000
$ws.Cells.Item(670, 16).Value2 = 17000
$ws.Cells.Item(670, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(670, 19).Value2 = 1417
$ws.Cells.Item(670, 20).Value2 = 12
$ws.Cells.Item(671, 4).Value2 = 44512
$ws.Cells.Item(671, 13).Value2 = 162
$ws.Cells.Item(671, 14).Value2 = 19000
$ws.Cells.Item(671, 15).Value2 = 19000
$ws.Cells.Item(671, 16).Value2 = 19000
$ws.Cells.Item(671, 19).Value2 = 1583
$ws.Cells.Item(672, 4).Value2 = 44512
$ws.Cells.Item(672, 12).Value2 = 'Segunda'
$ws.Cells.Item(672, 13).Value2 = 108
$ws.Cells.Item(672, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(672, 19).Value2 = 1357
$ws.Cells.Item(672, 20).Value2 = 14
$ws.Cells.Item(673, 4).Value2 = 44397
$ws.Cells.Item(673, 12).Value2 = 'Primera'
$ws.Cells.Item(673, 14).Value2 = 17000
$ws.Cells.Item(673, 15).Value2 = 17000
$ws.Cells.Item(673, 16).Value2 = 17000
$ws.Cells.Item(673, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(673, 19).Value2 = 1417
$ws.Cells.Item(673, 20).Value2 = 12
$ws.Cells.Item(674, 4).Value2 = 44181
$ws.Cells.Item(674, 14).Value2 = 19000
$ws.Cells.Item(674, 15).Value2 = 19000
$ws.Cells.Item(674, 16).Value2 = 19000
$ws.Cells.Item(674, 19).Value2 = 1583
$ws.Cells.Item(675, 4).Value2 = 44181
$ws.Cells.Item(675, 13).Value2 = 108
$ws.Cells.Item(675, 14).Value2 = 19000
$ws.Cells.Item(675, 15).Value2 = 19000
$ws.Cells.Item(675, 16).Value2 = 19000
$ws.Cells.Item(675, 19).Value2 = 1357
$ws.Cells.Item(676, 4).Value2 = 44285
$ws.Cells.Item(676, 13).Value2 = 162
$ws.Cells.Item(677, 4).Value2 = 44285
$ws.Cells.Item(677, 12).Value2 = 'Segunda'
$ws.Cells.Item(677, 13).Value2 = 54
$ws.Cells.Item(677, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(677, 19).Value2 = 1143
$ws.Cells.Item(677, 20).Value2 = 14
$ws.Cells.Item(678, 4).Value2 = 44362
$ws.Cells.Item(678, 13).Value2 = 108
$ws.Cells.Item(678, 19).Value2 = 1333
$ws.Cells.Item(679, 12).Value2 = 'Primera'
$ws.Cells.Item(679, 17).Value2 = '$/caja 10 unidades'
$ws.Cells.Item(679, 19).Value2 = 1700
$ws.Cells.Item(679, 20).Value2 = 10
$ws.Cells.Item(680, 4).Value2 = 44557
$ws.Cells.Item(680, 12).Value2 = 'Primera'
$ws.Cells.Item(680, 13).Value2 = 162
$ws.Cells.Item(680, 14).Value2 = 17000
$ws.Cells.Item(680, 15).Value2 = 17000
$ws.Cells.Item(680, 16).Value2 = 17000
$ws.Cells.Item(680, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(680, 19).Value2 = 1417
$ws.Cells.Item(680, 20).Value2 = 12
$ws.Cells.Item(681, 4).Value2 = 44557
$ws.Cells.Item(681, 12).Value2 = 'Segunda'
$ws.Cells.Item(681, 14).Value2 = 17000
$ws.Cells.Item(681, 15).Value2 = 17000
$ws.Cells.Item(681, 16).Value2 = 17000
$ws.Cells.Item(681, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(681, 19).Value2 = 1214
$ws.Cells.Item(681, 20).Value2 = 14
$ws.Cells.Item(682, 4).Value2 = 44747
$ws.Cells.Item(682, 12).Value2 = 'Especial'
$ws.Cells.Item(682, 13).Value2 = 54
$ws.Cells.Item(682, 14).Value2 = 22000
$ws.Cells.Item(682, 15).Value2 = 22000
$ws.Cells.Item(682, 16).Value2 = 22000
$ws.Cells.Item(682, 17).Value2 = '$/caja 10 unidades'
$ws.Cells.Item(682, 19).Value2 = 2200
$ws.Cells.Item(682, 20).Value2 = 10
$ws.Cells.Item(683, 4).Value2 = 44357
$ws.Cells.Item(683, 13).Value2 = 108
$ws.Cells.Item(683, 14).Value2 = 16000
$ws.Cells.Item(683, 15).Value2 = 16000
$ws.Cells.Item(683, 16).Value2 = 16000
$ws.Cells.Item(683, 19).Value2 = 1333
$ws.Cells.Item(684, 4).Value2 = 44357
$ws.Cells.Item(684, 13).Value2 = 108
$ws.Cells.Item(684, 14).Value2 = 16000
$ws.Cells.Item(684, 15).Value2 = 16000
$ws.Cells.Item(684, 16).Value2 = 16000
$ws.Cells.Item(684, 19).Value2 = 1143
$ws.Cells.Item(685, 4).Value2 = 44279
$ws.Cells.Item(685, 14).Value2 = 17000
$ws.Cells.Item(685, 15).Value2 = 17000
$ws.Cells.Item(685, 16).Value2 = 17000
$ws.Cells.Item(685, 19).Value2 = 1417
$ws.Cells.Item(686, 4).Value2 = 44279
$ws.Cells.Item(686, 14).Value2 = 17000
$ws.Cells.Item(686, 15).Value2 = 17000
$ws.Cells.Item(686, 16).Value2 = 17000
$ws.Cells.Item(686, 19).Value2 = 1214
$ws.Cells.Item(687, 4).Value2 = 44551
$ws.Cells.Item(687, 14).Value2 = 18000
$ws.Cells.Item(687, 15).Value2 = 18000
$ws.Cells.Item(687, 16).Value2 = 18000
$ws.Cells.Item(687, 19).Value2 = 1500
$ws.Cells.Item(688, 4).Value2 = 44551
$ws.Cells.Item(688, 13).Value2 = 162
$ws.Cells.Item(688, 14).Value2 = 18000
$ws.Cells.Item(688, 15).Value2 = 18000
$ws.Cells.Item(688, 16).Value2 = 18000
$ws.Cells.Item(688, 19).Value2 = 1286
$ws.Cells.Item(689, 4).Value2 = 44517
$ws.Cells.Item(689, 13).Value2 = 162
$ws.Cells.Item(689, 14).Value2 = 19000
$ws.Cells.Item(689, 15).Value2 = 19000
$ws.Cells.Item(689, 16).Value2 = 19000
$ws.Cells.Item(689, 19).Value2 = 1583
$ws.Cells.Item(690, 4).Value2 = 44517
$ws.Cells.Item(690, 14).Value2 = 19000
$ws.Cells.Item(690, 15).Value2 = 19000
$ws.Cells.Item(690, 16).Value2 = 19000
$ws.Cells.Item(690, 19).Value2 = 1357
$ws.Cells.Item(691, 4).Value2 = 44757
$ws.Cells.Item(691, 14).Value2 = 22000
$ws.Cells.Item(691, 15).Value2 = 22000
$ws.Cells.Item(691, 16).Value2 = 22000
$ws.Cells.Item(691, 19).Value2 = 1833
$ws.Cells.Item(692, 4).Value2 = 44757
$ws.Cells.Item(692, 14).Value2 = 22000
$ws.Cells.Item(692, 15).Value2 = 22000
$ws.Cells.Item(692, 16).Value2 = 22000
$ws.Cells.Item(692, 19).Value2 = 1571
$ws.Cells.Item(693, 4).Value2 = 44547
$ws.Cells.Item(693, 13).Value2 = 108
$ws.Cells.Item(693, 14).Value2 = 18000
$ws.Cells.Item(693, 15).Value2 = 18000
$ws.Cells.Item(693, 16).Value2 = 18000
$ws.Cells.Item(693, 19).Value2 = 1500
$ws.Cells.Item(694, 4).Value2 = 44547
$ws.Cells.Item(694, 12).Value2 = 'Segunda'
$ws.Cells.Item(694, 14).Value2 = 18000
$ws.Cells.Item(694, 15).Value2 = 18000
$ws.Cells.Item(694, 16).Value2 = 18000
$ws.Cells.Item(694, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(694, 19).Value2 = 1286
$ws.Cells.Item(694, 20).Value2 = 14
$ws.Cells.Item(695, 4).Value2 = 44321
$ws.Cells.Item(695, 12).Value2 = 'Primera'
$ws.Cells.Item(695, 13).Value2 = 54
$ws.Cells.Item(695, 14).Value2 = 16000
$ws.Cells.Item(695, 15).Value2 = 16000
$ws.Cells.Item(695, 16).Value2 = 16000
$ws.Cells.Item(695, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(695, 19).Value2 = 1333
$ws.Cells.Item(695, 20).Value2 = 12
$ws.Cells.Item(696, 12).Value2 = 'Primera'
$ws.Cells.Item(696, 13).Value2 = 108
$ws.Cells.Item(696, 14).Value2 = 22000
$ws.Cells.Item(696, 15).Value2 = 22000
$ws.Cells.Item(696, 16).Value2 = 22000
$ws.Cells.Item(696, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(696, 19).Value2 = 1833
$ws.Cells.Item(696, 20).Value2 = 12
$ws.Cells.Item(697, 4).Value2 = 44438
$ws.Cells.Item(697, 13).Value2 = 108
$ws.Cells.Item(697, 14).Value2 = 22000
$ws.Cells.Item(697, 15).Value2 = 22000
$ws.Cells.Item(697, 16).Value2 = 22000
$ws.Cells.Item(697, 19).Value2 = 1571
$ws.Cells.Item(698, 4).Value2 = 44438
$ws.Cells.Item(698, 12).Value2 = 'Segunda'
$ws.Cells.Item(698, 13).Value2 = 78
$ws.Cells.Item(698, 14).Value2 = 12000
$ws.Cells.Item(698, 15).Value2 = 12000
$ws.Cells.Item(698, 16).Value2 = 12000
$ws.Cells.Item(698, 17).Value2 = '$/caja 7 unidades'
$ws.Cells.Item(698, 18).Value2 = 'Ecuador'
$ws.Cells.Item(698, 19).Value2 = 1714
$ws.Cells.Item(698, 20).Value2 = 7
$ws.Cells.Item(699, 4).Value2 = 44355
$ws.Cells.Item(699, 18).Value2 = 'Ecuador'
$ws.Cells.Item(699, 19).Value2 = 1143
$ws.Cells.Item(700, 4).Value2 = 44657
$ws.Cells.Item(700, 14).Value2 = 18000
$ws.Cells.Item(700, 15).Value2 = 18000
$ws.Cells.Item(700, 16).Value2 = 18000
$ws.Cells.Item(700, 18).Value2 = 'Bolivia'
$ws.Cells.Item(700, 19).Value2 = 1500
$ws.Cells.Item(701, 4).Value2 = 44657
$ws.Cells.Item(701, 12).Value2 = 'Segunda'
$ws.Cells.Item(701, 13).Value2 = 54
$ws.Cells.Item(701, 14).Value2 = 18000
$ws.Cells.Item(701, 15).Value2 = 18000
$ws.Cells.Item(701, 16).Value2 = 18000
$ws.Cells.Item(701, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(701, 18).Value2 = 'Bolivia'
$ws.Cells.Item(701, 19).Value2 = 1286
$ws.Cells.Item(701, 20).Value2 = 14
$ws.Cells.Item(702, 4).Value2 = 44391
$ws.Cells.Item(702, 12).Value2 = 'Primera'
$ws.Cells.Item(702, 13).Value2 = 162
$ws.Cells.Item(702, 14).Value2 = 17000
$ws.Cells.Item(702, 15).Value2 = 17000
$ws.Cells.Item(702, 16).Value2 = 17000
$ws.Cells.Item(702, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(702, 19).Value2 = 1417
$ws.Cells.Item(702, 20).Value2 = 12
$ws.Cells.Item(703, 4).Value2 = 44453
$ws.Cells.Item(703, 14).Value2 = 22000
$ws.Cells.Item(703, 15).Value2 = 22000
$ws.Cells.Item(703, 16).Value2 = 22000
$ws.Cells.Item(703, 19).Value2 = 1833
$ws.Cells.Item(704, 4).Value2 = 44453
$ws.Cells.Item(704, 14).Value2 = 22000
$ws.Cells.Item(704, 15).Value2 = 22000
$ws.Cells.Item(704, 16).Value2 = 22000
$ws.Cells.Item(704, 19).Value2 = 1571
$ws.Cells.Item(705, 4).Value2 = 44186
$ws.Cells.Item(705, 14).Value2 = 18000
$ws.Cells.Item(705, 15).Value2 = 18000
$ws.Cells.Item(705, 16).Value2 = 18000
$ws.Cells.Item(705, 19).Value2 = 1500
$ws.Cells.Item(706, 4).Value2 = 44186
$ws.Cells.Item(706, 14).Value2 = 18000
$ws.Cells.Item(706, 15).Value2 = 18000
$ws.Cells.Item(706, 16).Value2 = 18000
$ws.Cells.Item(706, 19).Value2 = 1286
$ws.Cells.Item(707, 4).Value2 = 44189
$ws.Cells.Item(707, 12).Value2 = 'Primera'
$ws.Cells.Item(707, 14).Value2 = 17000
$ws.Cells.Item(707, 15).Value2 = 17000
$ws.Cells.Item(707, 16).Value2 = 17000
$ws.Cells.Item(707, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(707, 19).Value2 = 1417
$ws.Cells.Item(707, 20).Value2 = 12
$ws.Cells.Item(708, 4).Value2 = 44189
$ws.Cells.Item(708, 12).Value2 = 'Segunda'
$ws.Cells.Item(708, 13).Value2 = 108
$ws.Cells.Item(708, 14).Value2 = 17000
$ws.Cells.Item(708, 15).Value2 = 17000
$ws.Cells.Item(708, 16).Value2 = 17000
$ws.Cells.Item(708, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(708, 19).Value2 = 1214
$ws.Cells.Item(708, 20).Value2 = 14
$ws.Cells.Item(709, 12).Value2 = 'Especial'
$ws.Cells.Item(709, 13).Value2 = 108
$ws.Cells.Item(709, 17).Value2 = '$/caja 10 unidades'
$ws.Cells.Item(709, 19).Value2 = 1800
$ws.Cells.Item(709, 20).Value2 = 10
$ws.Cells.Item(710, 4).Value2 = 44609
$ws.Cells.Item(710, 12).Value2 = 'Primera'
$ws.Cells.Item(710, 13).Value2 = 162
$ws.Cells.Item(710, 14).Value2 = 18000
$ws.Cells.Item(710, 15).Value2 = 18000
$ws.Cells.Item(710, 16).Value2 = 18000
$ws.Cells.Item(710, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(710, 19).Value2 = 1500
$ws.Cells.Item(710, 20).Value2 = 12
$ws.Cells.Item(711, 4).Value2 = 44609
$ws.Cells.Item(711, 12).Value2 = 'Segunda'
$ws.Cells.Item(711, 14).Value2 = 18000
$ws.Cells.Item(711, 15).Value2 = 18000
$ws.Cells.Item(711, 16).Value2 = 18000
$ws.Cells.Item(711, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(711, 19).Value2 = 1286
$ws.Cells.Item(711, 20).Value2 = 14
$ws.Cells.Item(712, 12).Value2 = 'Especial'
$ws.Cells.Item(712, 17).Value2 = '$/caja 10 unidades'
$ws.Cells.Item(712, 19).Value2 = 2300
$ws.Cells.Item(712, 20).Value2 = 10
$ws.Cells.Item(713, 4).Value2 = 44489
$ws.Cells.Item(713, 12).Value2 = 'Primera'
$ws.Cells.Item(713, 13).Value2 = 162
$ws.Cells.Item(713, 14).Value2 = 23000
$ws.Cells.Item(713, 15).Value2 = 23000
$ws.Cells.Item(713, 16).Value2 = 23000
$ws.Cells.Item(713, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(713, 19).Value2 = 1917
$ws.Cells.Item(713, 20).Value2 = 12
$ws.Cells.Item(714, 4).Value2 = 44489
$ws.Cells.Item(714, 12).Value2 = 'Segunda'
$ws.Cells.Item(714, 14).Value2 = 23000
$ws.Cells.Item(714, 15).Value2 = 23000
$ws.Cells.Item(714, 16).Value2 = 23000
$ws.Cells.Item(714, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(714, 19).Value2 = 1643
$ws.Cells.Item(714, 20).Value2 = 14
$ws.Cells.Item(715, 4).Value2 = 44358
$ws.Cells.Item(715, 12).Value2 = 'Segunda'
$ws.Cells.Item(715, 14).Value2 = 16000
$ws.Cells.Item(715, 15).Value2 = 16000
$ws.Cells.Item(715, 16).Value2 = 16000
$ws.Cells.Item(715, 17).Value2 = '$/caja 14 unidades'
$ws.Cells.Item(715, 19).Value2 = 1143
$ws.Cells.Item(715, 20).Value2 = 14
$ws.Cells.Item(716, 1).Value2 = 3
$ws.Cells.Item(716, 2).Value2 = 'Femacal de La Calera'
$ws.Cells.Item(716, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(716, 4).Value2 = 44572
$ws.Cells.Item(716, 5).Value2 = 5
$ws.Cells.Item(716, 6).Value2 = 'Fruta'
$ws.Cells.Item(716, 7).Value2 = 100108
$ws.Cells.Item(716, 8).Value2 = 'Tropicales y subtropicales'
$ws.Cells.Item(716, 9).Value2 = 100108005
$ws.Cells.Item(716, 10).Value2 = 'Piña'
$ws.Cells.Item(716, 11).Value2 = 'Caramelo'
$ws.Cells.Item(716, 12).Value2 = 'Especial'
$ws.Cells.Item(716, 13).Value2 = 108
$ws.Cells.Item(716, 14).Value2 = 18000
$ws.Cells.Item(716, 15).Value2 = 18000
$ws.Cells.Item(716, 16).Value2 = 18000
$ws.Cells.Item(716, 17).Value2 = '$/caja 10 unidades'
$ws.Cells.Item(716, 18).Value2 = 'Ecuador'
$ws.Cells.Item(716, 19).Value2 = 1800
$ws.Cells.Item(716, 20).Value2 = 10
$ws.Cells.Item(717, 1).Value2 = 3
$ws.Cells.Item(717, 2).Value2 = 'Femacal de La Calera'
$ws.Cells.Item(717, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(717, 4).Value2 = 44572
$ws.Cells.Item(717, 5).Value2 = 5
$ws.Cells.Item(717, 6).Value2 = 'Fruta'
$ws.Cells.Item(717, 7).Value2 = 100108
$ws.Cells.Item(717, 8).Value2 = 'Tropicales y subtropicales'
$ws.Cells.Item(717, 9).Value2 = 100108005
$ws.Cells.Item(717, 10).Value2 = 'Piña'
$ws.Cells.Item(717, 11).Value2 = 'Caramelo'
$ws.Cells.Item(717, 12).Value2 = 'Primera'
$ws.Cells.Item(717, 13).Value2 = 108
$ws.Cells.Item(717, 14).Value2 = 18000
$ws.Cells.Item(717, 15).Value2 = 18000
$ws.Cells.Item(717, 16).Value2 = 18000
$ws.Cells.Item(717, 17).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(717, 18).Value2 = 'Ecuador'
$ws.Cells.Item(717, 19).Value2 = 1500
$ws.Cells.Item(717, 20).Value2 = 12
